# Updates the cryptos price/volume table (B/C/D/E columns, rows 2-51) with
# refreshed values. Price cells whose text looks like a plain number (e.g.
# "1.00", "0.999") are prefixed with a leading apostrophe so Excel stores
# them as literal text instead of silently converting them to a numeric
# value (which would lose the trailing zeros / exact formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.770.48"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.736.37"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'350.52"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").Value = "'106.82"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  -2.84%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.572"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "'19.33"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "3.165.73"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "2.742.09"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "50.765.30"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "'7.62"
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("D21").Value = "'12.79"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").Value = "0.0₃0947"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "'68.66"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "'261.44"
$ws.Range("E24").Value = "  -3.70%  "
$ws.Range("D25").Value = "'2.68"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'25.54"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").Value = "'0.158"
$ws.Range("E28").Value = "  +12.86%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'9.94"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").Value = "'51.65"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "'34.06"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'5.93"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("D34").Value = "'0.0434"
$ws.Range("E34").Value = "  -7.13%  "
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("D36").Value = "'5.16"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'18.38"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").Value = "'3.11"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").Value = "'2.43"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").Value = "'120.62"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Value = "'21.65"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "2.063.88"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.39"
$ws.Range("E49").Value = "  -5.97%  "
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "'0.900"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'57.51"
$ws.Range("E51").Value = "  -2.50%  "
